# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 532
$ws1.Range("G2").Value = 60
$ws1.Range("F3").Value = 6370
$ws1.Range("F5").Value = 97
$ws1.Range("F6").Value = 132
$ws1.Range("G10").Value = 138

# Sheet "全部类型" (all types, aggregated)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 532
$ws4.Range("G2").Value = 60
$ws4.Range("F3").Value = 6370
$ws4.Range("F6").Value = 97
$ws4.Range("F7").Value = 132
$ws4.Range("G12").Value = 138
